# Update countries table (Pais sheet) with refreshed COVID-19 data
# and re-sort two pairs of rows whose totals crossed over:
#   - Uruguay overtook Belice (rows 154/155)
#   - Montserrat overtook Islas Malvinas (rows 215/216)
# Also refreshes the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 01:20"

# --- Row 4 ---
$ws.Range("B4").Value = 7633570
$ws.Range("C4").Value = 31105
$ws.Range("D4").Value = 4843854
$ws.Range("E4").Value = 2575123
$ws.Range("G4").Value = 314
$ws.Range("H4").Value = 214593

# --- Row 5 ---
$ws.Range("B5").Value = 6622180
$ws.Range("C5").Value = 74767
$ws.Range("D5").Value = 5583453
$ws.Range("E5").Value = 936013

# --- Row 6 ---
$ws.Range("E6").Value = 505706
$ws.Range("G6").Value = 364
$ws.Range("H6").Value = 146375

# --- Row 8 ---
$ws.Range("B8").Value = 855052
$ws.Range("C8").Value = 6905
$ws.Range("D8").Value = 761674
$ws.Range("E8").Value = 66666
$ws.Range("G8").Value = 156
$ws.Range("H8").Value = 26712

# --- Row 9 ---
$ws.Range("B9").Value = 828169
$ws.Range("C9").Value = 3184
$ws.Range("D9").Value = 706223
$ws.Range("E9").Value = 89204
$ws.Range("G9").Value = 77
$ws.Range("H9").Value = 32742

# --- Row 11 ---
$ws.Range("B11").Value = 798486
$ws.Range("C11").Value = 7668
$ws.Range("D11").Value = 636672
$ws.Range("E11").Value = 140796
$ws.Range("G11").Value = 223
$ws.Range("H11").Value = 21018

# --- Row 15 ---
$ws.Range("B15").Value = 502978
$ws.Range("C15").Value = 7982
$ws.Range("G15").Value = 33
$ws.Range("H15").Value = 42350

# --- Row 37 ---
$ws.Range("B37").Value = 115286
$ws.Range("C37").Value = 633
$ws.Range("D37").Value = 91809
$ws.Range("E37").Value = 21054
$ws.Range("G37").Value = 9
$ws.Range("H37").Value = 2423

# --- Row 58 ---
$ws.Range("B58").Value = 59345
$ws.Range("C58").Value = 58
$ws.Range("D58").Value = 50768
$ws.Range("E58").Value = 7464

# --- Row 84 ---
$ws.Range("B84").Value = 21587
$ws.Range("C84").Value = 69
$ws.Range("D84").Value = 15014
$ws.Range("E84").Value = 5729
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 844

# --- Row 95 ---
$ws.Range("B95").Value = 14457
$ws.Range("C95").Value = 95
$ws.Range("E95").Value = 2992

# --- Row 115 ---
$ws.Range("B115").Value = 7520
$ws.Range("C115").Value = 3
$ws.Range("D115").Value = 7188
$ws.Range("E115").Value = 170

# --- Row 118 ---
$ws.Range("B118").Value = 6360
$ws.Range("C118").Value = 64
$ws.Range("D118").Value = 5416
$ws.Range("E118").Value = 879
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 65

# --- Row 135 ---
$ws.Range("B135").Value = 4366
$ws.Range("C135").Value = 37
$ws.Range("D135").Value = 1155
$ws.Range("E135").Value = 3006
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 205

# --- Row 154 ---
$ws.Range("A154").Value = "Uruguay"
$ws.Range("B154").Value = 2145
$ws.Range("C154").Value = 23
$ws.Range("D154").Value = 1844
$ws.Range("E154").Value = 253
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 48

# --- Row 155 ---
$ws.Range("A155").Value = "Belice"
$ws.Range("B155").Value = 2131
$ws.Range("C155").Value = 51
$ws.Range("D155").Value = 1346
$ws.Range("E155").Value = 756
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 29

# --- Row 215 ---
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# --- Row 216 ---
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
